$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Theory" row - pages written so far (B4): 12 -> 14
$ws.Range("B4").Value = 14

# Daily pages-written tracker, day 14 / 2021-09-24 (I14): 0 -> 2
$ws.Range("I14").Value = 2

# Move the active cell selection to I14, matching the saved view state
$ws.Range("I14").Select()

$wb.Save()
